$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("A")

# --- Row 5 (Lay1 Orbit block, index 15) ---
$ws.Range("B5").Value = 1.3654
$ws.Range("C5").Value = 257.9305
$ws.Range("E5").Value = "0.2540    0.2540    7.7724    0.0254    0.0254"
$ws.Range("L5").Value = 2.3042
$ws.Range("M5").Value = 435.2709
$ws.Range("O5").Value = "1.0922    1.0922    7.5692    0.0254    0.0254"

# --- Row 6 (Lay1 Orbit block, index 18) ---
$ws.Range("B6").Value = 1.2313
$ws.Range("C6").Value = 334.8874
$ws.Range("E6").Value = "0.2540    0.2540    6.5532    0.0254    0.0254"
$ws.Range("L6").Value = 1.5554
$ws.Range("M6").Value = 423.0375
$ws.Range("O6").Value = "0.3810    0.3810    7.5184    0.0254    0.0254"

# --- Row 9 (Lay2 Orbit block, index 6): clear computed columns ---
$ws.Range("B9").ClearContents() | Out-Null
$ws.Range("L9").ClearContents() | Out-Null

# --- Row 10 (Lay2 Orbit block, index 9): clear computed columns ---
$ws.Range("B10").ClearContents() | Out-Null
$ws.Range("L10").ClearContents() | Out-Null

# --- Row 11 (Lay2 Orbit block, index 12): clear computed columns ---
$ws.Range("B11").ClearContents() | Out-Null
$ws.Range("L11").ClearContents() | Out-Null

# --- Row 12 (Lay2 Orbit block, index 15): clear B/C/E, update L/M/O ---
$ws.Range("B12").ClearContents() | Out-Null
$ws.Range("C12").ClearContents() | Out-Null
$ws.Range("E12").ClearContents() | Out-Null
$ws.Range("L12").Value = 3.1617
$ws.Range("M12").Value = 597.2572
$ws.Range("O12").Value = "0.2540    0.2540   15.5956    0.0254    0.0254"

# --- Row 13 (Lay2 Orbit block, index 18): clear B/C/E, update L/M/O ---
$ws.Range("B13").ClearContents() | Out-Null
$ws.Range("C13").ClearContents() | Out-Null
$ws.Range("E13").ClearContents() | Out-Null
$ws.Range("L13").Value = 3.0753
$ws.Range("M13").Value = 836.431
$ws.Range("O13").Value = "0.2540    0.2540   15.0876    0.0254    0.0254"

# --- Update selection to M12 and make sure sheet "A" stays the active tab ---
$ws.Activate() | Out-Null
$ws.Range("M12").Select() | Out-Null
